$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the comment text for "Participants to exclude" (row 23) to mention MATLAB colon notation
$ws.Range("D23").Value = 'Separate with commas. Spaces will be ignored. Can use participant numbers (e.g., "1, 2") or IDs (e.g., "P1, P2" or "P01, P02" or "AB12, CD23"). Can use MATLAB''s colon notatation (e.g., "1, 3:10" excludes 1 and 3-through-10).'

# Update the comment text for "Runs to exclude in all participants (number)" (row 24) to mention MATLAB colon notation
$ws.Range("D24").Value = 'Separate with commas. Spaces will be ignored. Can use MATLAB''s colon notatation (e.g., "7, 10:12" excludes 7 and 10-through-12).'

# Adjust row heights to accommodate the new, longer text
$ws.Rows(23).RowHeight = 45
$ws.Rows(24).RowHeight = 30

# Update the view to reflect the scroll position / selection used when the edit was made
$ws.Application.Goto($ws.Range("C23"), $true)
$ws.Range("C23").Select()
$ws.Application.ActiveWindow.ScrollRow = 16
